$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update input values
$ws.Range("C2").Value = 17
$ws.Range("G2").Value = 1
$ws.Range("C3").Value = 1
$ws.Range("G3").Value = 10
$ws.Range("G4").Value = 3
$ws.Range("H4").Value = 0
$ws.Range("C5").Value = 2
$ws.Range("H5").Value = 2
$ws.Range("H6").Value = 0
$ws.Range("H7").Value = 2

# Update the selected cell/range to match the saved view state
$ws.Range("D11").Select()
